$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values below are crypto price/volume snapshots refreshed by the
# scheduled GitHub Actions scraper job. Price cells that look like plain
# numbers are written with a leading apostrophe (quote-prefix) so Excel
# stores them as literal text (preserving formats like "227.70" or
# "1.20" instead of collapsing to 227.7 / 1.2), matching the source feed.

$ws.Range("D2").Value = '38.809.91'
$ws.Range("E2").Value = '  +1.10%  '

$ws.Range("D3").Value = '2.106.22'
$ws.Range("E3").Value = '  +0.81%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '''227.70'

$ws.Range("D6").Value = '''0.616'
$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("D7").Value = '''62.36'
$ws.Range("E7").Value = '  +2.54%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  +2.25%  '

$ws.Range("D10").Value = '''0.0842'
$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("D12").Value = '''15.81'
$ws.Range("E12").Value = '  +6.51%  '

$ws.Range("D13").Value = '2.417.66'
$ws.Range("E13").Value = '  +0.89%  '

$ws.Range("D14").Value = '''22.01'
$ws.Range("E14").Value = '  -1.42%  '

$ws.Range("D15").Value = '''0.811'
$ws.Range("E15").Value = '  +3.44%  '

$ws.Range("E16").Value = '  +1.62%  '

$ws.Range("D17").Value = '2.104.36'
$ws.Range("E17").Value = '  +0.59%  '

$ws.Range("D18").Value = '38.823.20'
$ws.Range("E18").Value = '  +1.29%  '

$ws.Range("D19").Value = '''6.14'
$ws.Range("E19").Value = '  +1.34%  '

$ws.Range("D20").Value = '''71.62'
$ws.Range("E20").Value = '  +0.80%  '

$ws.Range("D21").Value = '0.0₃0844'
$ws.Range("E21").Value = '  +1.45%  '

$ws.Range("D22").Value = '''228.48'
$ws.Range("E22").Value = '  +1.43%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("E24").Value = '  -3.11%  '

$ws.Range("D25").Value = '''2.31'
$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''9.67'
$ws.Range("E26").Value = '  +2.67%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''172.26'
$ws.Range("E27").Value = '  +1.57%  '

$ws.Range("E28").Value = '  +2.04%  '

$ws.Range("D29").Value = '''1.43'
$ws.Range("E29").Value = '  +4.88%  '

$ws.Range("D30").Value = '''19.35'
$ws.Range("E30").Value = '  +1.74%  '

$ws.Range("D31").Value = '''2.56'
$ws.Range("E31").Value = '  +8.77%  '

$ws.Range("E32").Value = '  +0.59%  '

$ws.Range("B33").Value = 'THORChain'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D33").Value = '''7.20'
$ws.Range("E33").Value = '  +12.03%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''4.57'
$ws.Range("E34").Value = '  +1.27%  '

$ws.Range("E35").Value = '  -1.27%  '

$ws.Range("D36").Value = '''0.0618'
$ws.Range("E36").Value = '  +2.11%  '

$ws.Range("E37").Value = '  +0.25%  '

$ws.Range("D38").Value = '''3.54'
$ws.Range("E38").Value = '  -0.25%  '

$ws.Range("D39").Value = '''0.998'
$ws.Range("E39").Value = '  -0.22%  '

$ws.Range("D40").Value = '''18.17'
$ws.Range("E40").Value = '  -1.65%  '

$ws.Range("D41").Value = '''102.62'
$ws.Range("E41").Value = '  +2.60%  '

$ws.Range("E42").Value = '  +3.83%  '

$ws.Range("D43").Value = '1.528.61'
$ws.Range("E43").Value = '  -0.73%  '

$ws.Range("D44").Value = '''1.20'
$ws.Range("E44").Value = '  +7.51%  '

$ws.Range("E45").Value = '  -0.74%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '''0.0916'
$ws.Range("E46").Value = '  -1.82%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '''7.80'
$ws.Range("E47").Value = '  -0.03%  '

$ws.Range("E48").Value = '  +4.60%  '

$ws.Range("D49").Value = '''4.17'
$ws.Range("E49").Value = '  +0.73%  '

$ws.Range("E50").Value = '  -0.42%  '

$ws.Range("D51").Value = '2.304.98'
$ws.Range("E51").Value = '  +0.93%  '
